# Applies the "Natmi following Dr Hou advice" update:
# recomputed Efemp1->Egfr LR-pair statistics and an added "ECs" sending/target cluster,
# expanding the data block from 6 rows (2 clusters) to 9 rows (3 clusters).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old data (previously rows 2:7) and make room for the new 9-row block (rows 2:10).
$ws.Range("A2:T10").ClearContents()

# Populate the label columns (A-D) one full column at a time so the workbook's shared-string
# table is rebuilt in the same order as in the target file: ECs, FAPs, sCs, Efemp1, Egfr.

# Column A
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(10, 1).Value = "sCs"

# Column B
$ws.Cells.Item(2, 2).Value = "Efemp1"
$ws.Cells.Item(3, 2).Value = "Efemp1"
$ws.Cells.Item(4, 2).Value = "Efemp1"
$ws.Cells.Item(5, 2).Value = "Efemp1"
$ws.Cells.Item(6, 2).Value = "Efemp1"
$ws.Cells.Item(7, 2).Value = "Efemp1"
$ws.Cells.Item(8, 2).Value = "Efemp1"
$ws.Cells.Item(9, 2).Value = "Efemp1"
$ws.Cells.Item(10, 2).Value = "Efemp1"

# Column C
$ws.Cells.Item(2, 3).Value = "Egfr"
$ws.Cells.Item(3, 3).Value = "Egfr"
$ws.Cells.Item(4, 3).Value = "Egfr"
$ws.Cells.Item(5, 3).Value = "Egfr"
$ws.Cells.Item(6, 3).Value = "Egfr"
$ws.Cells.Item(7, 3).Value = "Egfr"
$ws.Cells.Item(8, 3).Value = "Egfr"
$ws.Cells.Item(9, 3).Value = "Egfr"
$ws.Cells.Item(10, 3).Value = "Egfr"

# Column D
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(10, 4).Value = "sCs"

# Populate the numeric statistic columns (E-T).

# Column E
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(10, 5).Value = 3

# Column F
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(10, 6).Value = 1

# Column G
$ws.Cells.Item(2, 7).Value = 1.181524666666667
$ws.Cells.Item(3, 7).Value = 1.181524666666667
$ws.Cells.Item(4, 7).Value = 1.181524666666667
$ws.Cells.Item(5, 7).Value = 168.5040536666667
$ws.Cells.Item(6, 7).Value = 168.5040536666667
$ws.Cells.Item(7, 7).Value = 168.5040536666667
$ws.Cells.Item(8, 7).Value = 4.891206
$ws.Cells.Item(9, 7).Value = 4.891206
$ws.Cells.Item(10, 7).Value = 4.891206

# Column H
$ws.Cells.Item(2, 8).Value = 3.544574
$ws.Cells.Item(3, 8).Value = 3.544574
$ws.Cells.Item(4, 8).Value = 3.544574
$ws.Cells.Item(5, 8).Value = 505.512161
$ws.Cells.Item(6, 8).Value = 505.512161
$ws.Cells.Item(7, 8).Value = 505.512161
$ws.Cells.Item(8, 8).Value = 14.673618
$ws.Cells.Item(9, 8).Value = 14.673618
$ws.Cells.Item(10, 8).Value = 14.673618

# Column I
$ws.Cells.Item(2, 9).Value = 0.006767936934905889
$ws.Cells.Item(3, 9).Value = 0.006767936934905889
$ws.Cells.Item(4, 9).Value = 0.006767936934905889
$ws.Cells.Item(5, 9).Value = 0.9652145576523421
$ws.Cells.Item(6, 9).Value = 0.9652145576523421
$ws.Cells.Item(7, 9).Value = 0.9652145576523421
$ws.Cells.Item(8, 9).Value = 0.02801750541275197
$ws.Cells.Item(9, 9).Value = 0.02801750541275197
$ws.Cells.Item(10, 9).Value = 0.02801750541275197

# Column J
$ws.Cells.Item(2, 10).Value = 0.00676793693490589
$ws.Cells.Item(3, 10).Value = 0.00676793693490589
$ws.Cells.Item(4, 10).Value = 0.00676793693490589
$ws.Cells.Item(5, 10).Value = 0.9652145576523421
$ws.Cells.Item(6, 10).Value = 0.9652145576523421
$ws.Cells.Item(7, 10).Value = 0.9652145576523421
$ws.Cells.Item(8, 10).Value = 0.02801750541275197
$ws.Cells.Item(9, 10).Value = 0.02801750541275197
$ws.Cells.Item(10, 10).Value = 0.02801750541275197

# Column K
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(10, 11).Value = 3

# Column L
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(10, 12).Value = 1

# Column M
$ws.Cells.Item(2, 13).Value = 2.08532
$ws.Cells.Item(3, 13).Value = 101.898173
$ws.Cells.Item(4, 13).Value = 14.921347
$ws.Cells.Item(5, 13).Value = 2.08532
$ws.Cells.Item(6, 13).Value = 101.898173
$ws.Cells.Item(7, 13).Value = 14.921347
$ws.Cells.Item(8, 13).Value = 2.08532
$ws.Cells.Item(9, 13).Value = 101.898173
$ws.Cells.Item(10, 13).Value = 14.921347

# Column N
$ws.Cells.Item(2, 14).Value = 6.25596
$ws.Cells.Item(3, 14).Value = 305.694519
$ws.Cells.Item(4, 14).Value = 44.76404100000001
$ws.Cells.Item(5, 14).Value = 6.25596
$ws.Cells.Item(6, 14).Value = 305.694519
$ws.Cells.Item(7, 14).Value = 44.76404100000001
$ws.Cells.Item(8, 14).Value = 6.25596
$ws.Cells.Item(9, 14).Value = 305.694519
$ws.Cells.Item(10, 14).Value = 44.76404100000001

# Column O
$ws.Cells.Item(2, 15).Value = 0.01753772176136817
$ws.Cells.Item(3, 15).Value = 0.8569724579756384
$ws.Cells.Item(4, 15).Value = 0.1254898202629935
$ws.Cells.Item(5, 15).Value = 0.01753772176136817
$ws.Cells.Item(6, 15).Value = 0.8569724579756384
$ws.Cells.Item(7, 15).Value = 0.1254898202629935
$ws.Cells.Item(8, 15).Value = 0.01753772176136817
$ws.Cells.Item(9, 15).Value = 0.8569724579756384
$ws.Cells.Item(10, 15).Value = 0.1254898202629935

# Column P
$ws.Cells.Item(2, 16).Value = 0.01753772176136816
$ws.Cells.Item(3, 16).Value = 0.8569724579756383
$ws.Cells.Item(4, 16).Value = 0.1254898202629935
$ws.Cells.Item(5, 16).Value = 0.01753772176136816
$ws.Cells.Item(6, 16).Value = 0.8569724579756383
$ws.Cells.Item(7, 16).Value = 0.1254898202629935
$ws.Cells.Item(8, 16).Value = 0.01753772176136816
$ws.Cells.Item(9, 16).Value = 0.8569724579756383
$ws.Cells.Item(10, 16).Value = 0.1254898202629935

# Column Q
$ws.Cells.Item(2, 17).Value = 2.463857017893333
$ws.Cells.Item(3, 17).Value = 120.3952048877673
$ws.Cells.Item(4, 17).Value = 17.62993954039267
$ws.Cells.Item(5, 17).Value = 351.3848731921733
$ws.Cells.Item(6, 17).Value = 17170.25521172728
$ws.Cells.Item(7, 17).Value = 2514.307455666956
$ws.Cells.Item(8, 17).Value = 10.19972969592
$ws.Cells.Item(9, 17).Value = 498.404955166638
$ws.Cells.Item(10, 17).Value = 72.98338197448201

# Column R
$ws.Cells.Item(2, 18).Value = 22.17471316104
$ws.Cells.Item(3, 18).Value = 1083.556843989906
$ws.Cells.Item(4, 18).Value = 158.669455863534
$ws.Cells.Item(5, 18).Value = 3162.46385872956
$ws.Cells.Item(6, 18).Value = 154532.2969055456
$ws.Cells.Item(7, 18).Value = 22628.7671010026
$ws.Cells.Item(8, 18).Value = 91.79756726328
$ws.Cells.Item(9, 18).Value = 4485.644596499743
$ws.Cells.Item(10, 18).Value = 656.8504377703381

# Column S
$ws.Cells.Item(2, 19).Value = 0.0001186941948628664
$ws.Cells.Item(3, 19).Value = 0.005799935550530408
$ws.Cells.Item(4, 19).Value = 0.0008493071895126153
$ws.Cells.Item(5, 19).Value = 0.01692766435212883
$ws.Cells.Item(6, 19).Value = 0.8271622919451961
$ws.Cells.Item(7, 19).Value = 0.1211246013550172
$ws.Cells.Item(8, 19).Value = 0.0004913632143764707
$ws.Cells.Item(9, 19).Value = 0.02401023047991181
$ws.Cells.Item(10, 19).Value = 0.003515911718463693

# Column T
$ws.Cells.Item(2, 20).Value = 0.0001186941948628664
$ws.Cells.Item(3, 20).Value = 0.005799935550530408
$ws.Cells.Item(4, 20).Value = 0.0008493071895126152
$ws.Cells.Item(5, 20).Value = 0.01692766435212882
$ws.Cells.Item(6, 20).Value = 0.827162291945196
$ws.Cells.Item(7, 20).Value = 0.1211246013550172
$ws.Cells.Item(8, 20).Value = 0.0004913632143764704
$ws.Cells.Item(9, 20).Value = 0.02401023047991181
$ws.Cells.Item(10, 20).Value = 0.003515911718463692
